$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove column A ("Id informe") - no longer needed; everything shifts left.
$ws.Columns.Item(1).Delete()

# 2. Rename the last column's header ("Proceso terminado?" -> "Estado proceso")
#    After the delete above, that column is now column O (15).
$ws.Range("O1").Value2 = "Estado proceso"

# 3. Move that "Estado proceso" column so it becomes the second column (right after "Id proceso").
$ws.Columns.Item(15).Cut()
$ws.Columns.Item(2).Insert()

# 4. Refresh the AutoFilter so it spans the new data range A1:O1.
$ws.AutoFilterMode = $false
$ws.Range("A1:O1").AutoFilter()

# 5. Update the workbook-level hidden defined name backing the filter (_FilterDatabase)
#    so it matches the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Data!`$A`$1:`$P`$1"
    }
}
